$wb = $excel.ActiveWorkbook

# Sheet "建物" (building) - fix property_category column (I) for data rows 2-12
# from "land" to "building"
$wsBuilding = $wb.Worksheets.Item("建物")
for ($r = 2; $r -le 12; $r++) {
    $wsBuilding.Range("I" + $r).Value = "building"
}

# Sheet "汽車" (car) - fix property_category column (H) for data rows 2-3
# from "land" to "car"
$wsCar = $wb.Worksheets.Item("汽車")
for ($r = 2; $r -le 3; $r++) {
    $wsCar.Range("H" + $r).Value = "car"
}
